$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 140.28572
$ws.Cells.Item(33, 9).Value = 140.28572
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 140.28572
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 88.71428

$ws.Cells.Item(42, 8).Value = 204
$ws.Cells.Item(42, 9).Value = 135.6
$ws.Cells.Item(42, 10).Value = 375
$ws.Cells.Item(42, 11).Value = 406.8
$ws.Cells.Item(42, 12).Value = 1125
$ws.Cells.Item(42, 13).Value = -176.8
$ws.Cells.Item(42, 14).Value = -1585

$ws.Cells.Item(96, 8).Value = 317.16666
$ws.Cells.Item(96, 9).Value = 209.63637
$ws.Cells.Item(96, 10).Value = 1500
$ws.Cells.Item(96, 11).Value = 628.9091100000001
$ws.Cells.Item(96, 12).Value = 4500
$ws.Cells.Item(96, 13).Value = 744.0908899999999
$ws.Cells.Item(96, 14).Value = -7246

$ws.Cells.Item(116, 8).Value = 3700
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 3700
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 3700
$ws.Cells.Item(116, 13).Value = ""
$ws.Cells.Item(116, 14).Value = -10584

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5104.15
$ws.Cells.Item(2, 9).Value = 3318.6875
$ws.Cells.Item(2, 10).Value = 12246
$ws.Cells.Item(2, 11).Value = 3318.6875
$ws.Cells.Item(2, 12).Value = 12246
$ws.Cells.Item(2, 13).Value = -3205.6875
$ws.Cells.Item(2, 14).Value = -12472

$ws.Cells.Item(32, 8).Value = 2947914.8
$ws.Cells.Item(32, 9).Value = 3689.074
$ws.Cells.Item(32, 10).Value = 14304214
$ws.Cells.Item(32, 11).Value = 3689.074
$ws.Cells.Item(32, 12).Value = 14304214
$ws.Cells.Item(32, 13).Value = -3402.074
$ws.Cells.Item(32, 14).Value = -14304788

$ws.Cells.Item(44, 8).Value = 11858.842
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 11858.842
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 11858.842
$ws.Cells.Item(44, 14).Value = -12834.842

$ws.Cells.Item(62, 8).Value = 120000
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 120000
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 120000
$ws.Cells.Item(62, 14).Value = -121248

$ws.Cells.Item(65, 8).Value = 120000
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 120000
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 360000
$ws.Cells.Item(65, 14).Value = -366240

$ws.Cells.Item(74, 8).Value = 4083.5
$ws.Cells.Item(74, 9).Value = 3378.0527
$ws.Cells.Item(74, 10).Value = 8551.333000000001
$ws.Cells.Item(74, 11).Value = 3378.0527
$ws.Cells.Item(74, 12).Value = 8551.333000000001
$ws.Cells.Item(74, 13).Value = -2504.0527
$ws.Cells.Item(74, 14).Value = -10299.333

$ws.Cells.Item(77, 8).Value = 4083.5
$ws.Cells.Item(77, 9).Value = 3378.0527
$ws.Cells.Item(77, 10).Value = 8551.333000000001
$ws.Cells.Item(77, 11).Value = 16890.2635
$ws.Cells.Item(77, 12).Value = 42756.665
$ws.Cells.Item(77, 13).Value = -12522.2635
$ws.Cells.Item(77, 14).Value = -51492.665

$ws.Cells.Item(116, 8).Value = 5104.15
$ws.Cells.Item(116, 9).Value = 3318.6875
$ws.Cells.Item(116, 10).Value = 12246
$ws.Cells.Item(116, 11).Value = 3318.6875
$ws.Cells.Item(116, 12).Value = 12246
$ws.Cells.Item(116, 13).Value = -1024.6875
$ws.Cells.Item(116, 14).Value = -16834

$ws.Cells.Item(132, 8).Value = 3558.875
$ws.Cells.Item(132, 9).Value = 2078.6667
$ws.Cells.Item(132, 10).Value = 7999.5
$ws.Cells.Item(132, 11).Value = 6236.000100000001
$ws.Cells.Item(132, 12).Value = 23998.5
$ws.Cells.Item(132, 13).Value = -3706.000100000001
$ws.Cells.Item(132, 14).Value = -29058.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5104.15
$ws.Cells.Item(3, 9).Value = 3318.6875
$ws.Cells.Item(3, 10).Value = 12246
$ws.Cells.Item(3, 11).Value = 3318.6875
$ws.Cells.Item(3, 12).Value = 12246
$ws.Cells.Item(3, 13).Value = -3204.6875
$ws.Cells.Item(3, 14).Value = -12474

$ws.Cells.Item(20, 8).Value = 2506.7144
$ws.Cells.Item(20, 9).Value = 2633.375
$ws.Cells.Item(20, 10).Value = 2101.4
$ws.Cells.Item(20, 11).Value = 2633.375
$ws.Cells.Item(20, 12).Value = 2101.4
$ws.Cells.Item(20, 13).Value = -2386.375
$ws.Cells.Item(20, 14).Value = -2595.4

$ws.Cells.Item(105, 8).Value = 1744.4615
$ws.Cells.Item(105, 9).Value = 1559.8889
$ws.Cells.Item(105, 10).Value = 2159.75
$ws.Cells.Item(105, 11).Value = 1559.8889
$ws.Cells.Item(105, 12).Value = 2159.75
$ws.Cells.Item(105, 13).Value = 187.1111000000001
$ws.Cells.Item(105, 14).Value = -5653.75

$ws.Cells.Item(134, 8).Value = 3300.8
$ws.Cells.Item(134, 9).Value = 2854.7693
$ws.Cells.Item(134, 10).Value = 6200
$ws.Cells.Item(134, 11).Value = 8564.3079
$ws.Cells.Item(134, 12).Value = 18600
$ws.Cells.Item(134, 13).Value = -6029.3079
$ws.Cells.Item(134, 14).Value = -23670

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4760
$ws.Cells.Item(58, 9).Value = 1219
$ws.Cells.Item(58, 10).Value = 7592.8
$ws.Cells.Item(58, 11).Value = 1219
$ws.Cells.Item(58, 12).Value = 7592.8
$ws.Cells.Item(58, 13).Value = -1016
$ws.Cells.Item(58, 14).Value = -7998.8

$ws.Cells.Item(69, 8).Value = 7500
$ws.Cells.Item(69, 9).Value = 7500
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 7500
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = -6751
$ws.Cells.Item(69, 14).Value = ""

$ws.Cells.Item(72, 8).Value = 7500
$ws.Cells.Item(72, 9).Value = 7500
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 22500
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = -18756
$ws.Cells.Item(72, 14).Value = ""

$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = ""

$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = ""

$ws.Cells.Item(99, 8).Value = 4204.5835
$ws.Cells.Item(99, 9).Value = 3938
$ws.Cells.Item(99, 10).Value = 4577.8
$ws.Cells.Item(99, 11).Value = 3938
$ws.Cells.Item(99, 12).Value = 4577.8
$ws.Cells.Item(99, 13).Value = -2440
$ws.Cells.Item(99, 14).Value = -7573.8

$ws.Cells.Item(126, 8).Value = 4204.5835
$ws.Cells.Item(126, 9).Value = 3938
$ws.Cells.Item(126, 10).Value = 4577.8
$ws.Cells.Item(126, 11).Value = 11814
$ws.Cells.Item(126, 12).Value = 13733.4
$ws.Cells.Item(126, 13).Value = -9344
$ws.Cells.Item(126, 14).Value = -18673.4

$ws.Cells.Item(136, 8).Value = 4760
$ws.Cells.Item(136, 9).Value = 1219
$ws.Cells.Item(136, 10).Value = 7592.8
$ws.Cells.Item(136, 11).Value = 3657
$ws.Cells.Item(136, 12).Value = 22778.4
$ws.Cells.Item(136, 13).Value = -1107
$ws.Cells.Item(136, 14).Value = -27878.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 4650
$ws.Cells.Item(137, 9).Value = 3000
$ws.Cells.Item(137, 10).Value = 6300
$ws.Cells.Item(137, 11).Value = 9000
$ws.Cells.Item(137, 12).Value = 18900
$ws.Cells.Item(137, 13).Value = -3900
$ws.Cells.Item(137, 14).Value = -29100

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 2536262
$ws.Cells.Item(34, 9).Value = 5012500
$ws.Cells.Item(34, 10).Value = 60024
$ws.Cells.Item(34, 11).Value = 5012500
$ws.Cells.Item(34, 12).Value = 60024
$ws.Cells.Item(34, 13).Value = -5012328
$ws.Cells.Item(34, 14).Value = -60368

$ws.Cells.Item(46, 8).Value = 6557.6875
$ws.Cells.Item(46, 9).Value = 6385.6
$ws.Cells.Item(46, 10).Value = 6635.909
$ws.Cells.Item(46, 11).Value = 6385.6
$ws.Cells.Item(46, 12).Value = 6635.909
$ws.Cells.Item(46, 13).Value = -6197.6
$ws.Cells.Item(46, 14).Value = -7011.909

$ws.Cells.Item(132, 8).Value = 3693.5
$ws.Cells.Item(132, 9).Value = 3817.125
$ws.Cells.Item(132, 10).Value = 3199
$ws.Cells.Item(132, 11).Value = 11451.375
$ws.Cells.Item(132, 12).Value = 9597
$ws.Cells.Item(132, 13).Value = -8921.375
$ws.Cells.Item(132, 14).Value = -14657

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 1420
$ws.Cells.Item(19, 9).Value = 1400
$ws.Cells.Item(19, 10).Value = 1500
$ws.Cells.Item(19, 11).Value = 1400
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = -1226
$ws.Cells.Item(19, 14).Value = -1848

$ws.Cells.Item(113, 8).Value = 762.0714
$ws.Cells.Item(113, 9).Value = 779
$ws.Cells.Item(113, 10).Value = 700
$ws.Cells.Item(113, 11).Value = 2337
$ws.Cells.Item(113, 12).Value = 2100
$ws.Cells.Item(113, 13).Value = -167
$ws.Cells.Item(113, 14).Value = -6440
